$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 24,3
$arr[0,0] = 1.038240949722478
$arr[0,1] = 0.2665211918717034
$arr[0,2] = 0.02380986417120567
$arr[1,0] = 0.944975639967879
$arr[1,1] = 0.2450261732752779
$arr[1,2] = 0.02346685714776697
$arr[2,0] = 0.8879779308421121
$arr[2,1] = 0.2317185442680341
$arr[2,2] = 0.0232559817544491
$arr[3,0] = 0.8648189441159673
$arr[3,1] = 0.2262680288233128
$arr[3,2] = 0.0231699898031259
$arr[4,0] = 0.8609775441794056
$arr[4,1] = 0.2253613114005475
$arr[4,2] = 0.02315570763072827
$arr[5,0] = 0.8876653241347583
$arr[5,1] = 0.2316451482361117
$arr[5,2] = 0.02325482226205722
$arr[6,0] = 1.006027961936979
$arr[6,1] = 0.2591324578286276
$arr[6,2] = 0.02369165604297407
$arr[7,0] = 1.240238733983801
$arr[7,1] = 0.3121673007282482
$arr[7,2] = 0.02454582167461794
$arr[8,0] = 1.413584839774728
$arr[8,1] = 0.3506092697668635
$arr[8,2] = 0.02517148620862031
$arr[9,0] = 1.49271964564133
$arr[9,1] = 0.3679855639208824
$arr[9,2] = 0.02545563034049891
$arr[10,0] = 1.522725541953264
$arr[10,1] = 0.3745495749961094
$arr[10,2] = 0.02556315265163178
$arr[11,0] = 1.516261502203179
$arr[11,1] = 0.3731366105013194
$arr[11,2] = 0.02553999935677354
$arr[12,0] = 1.495187466625566
$arr[12,1] = 0.3685259101914369
$arr[12,2] = 0.02546447785633177
$arr[13,0] = 1.482284100040772
$arr[13,1] = 0.365699636657439
$arr[13,2] = 0.02541820848432863
$arr[14,0] = 1.408418733293672
$arr[14,1] = 0.3494714487771944
$arr[14,2] = 0.02515290643421153
$arr[15,0] = 1.363175538603116
$arr[15,1] = 0.3394874905648635
$arr[15,2] = 0.02499002458195321
$arr[16,0] = 1.337179140095259
$arr[16,1] = 0.3337345217053667
$arr[16,2] = 0.02489629516572478
$arr[17,0] = 1.328381746940693
$arr[17,1] = 0.331784869710475
$arr[17,2] = 0.02486455269477617
$arr[18,0] = 1.367989038961753
$arr[18,1] = 0.3405513827767948
$arr[18,2] = 0.02500736825852812
$arr[19,0] = 1.501376362704491
$arr[19,1] = 0.3698806195095869
$arr[19,2] = 0.02548666250815046
$arr[20,0] = 1.588781162848989
$arr[20,1] = 0.388955612891408
$arr[20,2] = 0.0257994567757045
$arr[21,0] = 1.542110942768375
$arr[21,1] = 0.3787834823656908
$arr[21,2] = 0.0256325568301925
$arr[22,0] = 1.365812810170155
$arr[22,1] = 0.3400704377866361
$arr[22,2] = 0.02499952745001366
$arr[23,0] = 1.17665453309786
$arr[23,1] = 0.2979122197349398
$arr[23,2] = 0.02431505175663062
$ws.Range("B2:D25").Value = $arr

$arr = New-Object "object[,]" 24,4
$arr[0,0] = 0.9144269611472566
$arr[0,1] = 0.7683326004949862
$arr[0,2] = 0.8140152192382573
$arr[0,3] = 0.7720220758624876
$arr[1,0] = 0.9013263571936676
$arr[1,1] = 0.7556657309808656
$arr[1,2] = 0.8133260020227056
$arr[1,3] = 0.7748831900374213
$arr[2,0] = 0.8939441749435417
$arr[2,1] = 0.7485020716354711
$arr[2,2] = 0.8133783259870597
$arr[2,3] = 0.7771319533533756
$arr[3,0] = 0.8911018798137036
$arr[3,1] = 0.7457366757248138
$arr[3,2] = 0.81351910408317
$arr[3,3] = 0.7781718646007008
$arr[4,0] = 0.8906399342984486
$arr[4,1] = 0.7452867614949952
$arr[4,2] = 0.813549691698924
$arr[4,3] = 0.7783519955890696
$arr[5,0] = 0.8939051711070221
$arr[5,1] = 0.7484641542616544
$arr[5,2] = 0.8133797410528132
$arr[5,3] = 0.7771454780953846
$arr[6,0] = 0.9097723188588844
$arr[6,1] = 0.7638373496479289
$arr[6,2] = 0.8136788142414844
$arr[6,3] = 0.7729063442585584
$arr[7,0] = 0.9461573484679917
$arr[7,1] = 0.7988802476060215
$arr[7,2] = 0.8180446158427515
$arr[7,3] = 0.7685073134773148
$arr[8,0] = 0.9761336646762686
$arr[8,1] = 0.827650302995238
$arr[8,2] = 0.8235667856652071
$arr[8,3] = 0.7676761197814628
$arr[9,0] = 0.9904822378773304
$arr[9,1] = 0.8414038843660592
$arr[9,2] = 0.8265839911722139
$arr[9,3] = 0.7678225613394574
$arr[10,0] = 0.9960185697104862
$arr[10,1] = 0.8467083917800267
$arr[10,2] = 0.827799334705702
$arr[10,3] = 0.7679536955148691
$arr[11,0] = 0.9948216410317485
$arr[11,1] = 0.8455616781517961
$arr[11,2] = 0.8275343487560178
$arr[11,3] = 0.7679220835545735
$arr[12,0] = 0.9909356517503198
$arr[12,1] = 0.8418383556635263
$arr[12,2] = 0.8266825184412596
$arr[12,3] = 0.7678318316561317
$arr[13,0] = 0.9885687763733841
$arr[13,1] = 0.8395702741802324
$arr[13,2] = 0.8261702322161852
$arr[13,3] = 0.7677864129443321
$arr[14,0] = 0.9752103192452921
$arr[14,1] = 0.8267649226436617
$arr[14,2] = 0.8233797821191615
$arr[14,3] = 0.7676771261850774
$arr[15,0] = 0.9671980516108931
$arr[15,1] = 0.8190802177813197
$arr[15,2] = 0.8217974199774289
$arr[15,3] = 0.7677446107900252
$arr[16,0] = 0.9626566169910546
$arr[16,1] = 0.814722826372261
$arr[16,2] = 0.8209348242268675
$arr[16,3] = 0.767832784832251
$arr[17,0] = 0.9611304613132745
$arr[17,1] = 0.813258231142413
$arr[17,2] = 0.8206509239931705
$arr[17,3] = 0.7678711088689738
$arr[18,0] = 0.9680440328329638
$arr[18,1] = 0.8198917798635819
$arr[18,2] = 0.821960944098862
$arr[18,3] = 0.7677323166931345
$arr[19,0] = 0.9920742665942441
$arr[19,1] = 0.8429293670149036
$arr[19,2] = 0.826930744827024
$arr[19,3] = 0.7678562849014696
$arr[20,0] = 1.008379033579018
$arr[20,1] = 0.8585474778577407
$arr[20,2] = 0.8306031644994505
$arr[20,3] = 0.7683785423682608
$arr[21,0] = 0.9996218664388721
$arr[21,1] = 0.8501602147161691
$arr[21,2] = 0.828604243992686
$arr[21,3] = 0.7680593478237583
$arr[22,0] = 0.9676613625148747
$arr[22,1] = 0.8195246837732952
$arr[22,2] = 0.8218868680557989
$arr[22,3] = 0.7677377210642078
$arr[23,0] = 0.9357468096249164
$arr[23,1] = 0.7888720009531056
$arr[23,2] = 0.8164578539402925
$arr[23,3] = 0.769276836859774
$ws.Range("F2:I25").Value = $arr

$arr = New-Object "object[,]" 24,3
$arr[0,0] = 0.2688695164919324
$arr[0,1] = 0.2490549816902998
$arr[0,2] = 1.351258787217532
$arr[1,0] = 0.2670700693394394
$arr[1,1] = 0.2349287805085325
$arr[1,2] = 1.368168211889422
$arr[2,0] = 0.2661002671769239
$arr[2,1] = 0.2263614891806895
$arr[2,2] = 1.379093728034885
$arr[3,0] = 0.26573905805013
$arr[3,1] = 0.2228970803747288
$arr[3,2] = 1.383682451433669
$arr[4,0] = 0.2656811335396156
$arr[4,1] = 0.2223234418960089
$arr[4,2] = 1.384452646947256
$arr[5,0] = 0.2660952581071285
$arr[5,1] = 0.2263146581589481
$arr[5,2] = 1.379155060628467
$arr[6,0] = 0.2682210507381129
$arr[6,1] = 0.2441622777525652
$arr[6,2] = 1.356976403222683
$arr[7,0] = 0.2734606689013859
$arr[7,1] = 0.2800015227454935
$arr[7,2] = 1.317795658678744
$arr[8,0] = 0.2779629397373498
$arr[8,1] = 0.3068438203301582
$arr[8,2] = 1.291640745218952
$arr[9,0] = 0.2801529014015358
$arr[9,1] = 0.3191660794785918
$arr[9,2] = 1.280314176477027
$arr[10,0] = 0.2810025677362376
$arr[10,1] = 0.3238481745245991
$arr[10,2] = 1.276107366668235
$arr[11,0] = 0.2808186710569061
$arr[11,1] = 0.3228390953039266
$arr[11,2] = 1.277009715050347
$arr[12,0] = 0.2802223957418022
$arr[12,1] = 0.3195509598195372
$arr[12,2] = 1.279966429122585
$arr[13,0] = 0.2798598126824743
$arr[13,1] = 0.3175389519888228
$arr[13,2] = 1.281788226206306
$arr[14,0] = 0.2778226730959261
$arr[14,1] = 0.3060407664008764
$arr[14,2] = 1.292392474063181
$arr[15,0] = 0.2766092699731075
$arr[15,1] = 0.299015490561743
$arr[15,2] = 1.299044295647239
$arr[16,0] = 0.2759247061328978
$arr[16,1] = 0.2949852456871795
$arr[16,2] = 1.302924024456892
$arr[17,0] = 0.2756952187798305
$arr[17,1] = 0.2936224824389129
$arr[17,2] = 1.304246865892537
$arr[18,0] = 0.2767370569356302
$arr[18,1] = 0.2997622565020137
$arr[18,2] = 1.29833063113179
$arr[19,0] = 0.2803969834581039
$arr[19,1] = 0.3205163337232122
$arr[19,2] = 1.279095735319626
$arr[20,0] = 0.282907709722366
$arr[20,1] = 0.3341730790075061
$arr[20,2] = 1.267004489868665
$arr[21,0] = 0.281556829613308
$arr[21,1] = 0.3268757691912967
$arr[21,2] = 1.273413857880609
$arr[22,0] = 0.2766792438380179
$arr[22,1] = 0.2994246165752017
$arr[22,2] = 1.298653105957129
$arr[23,0] = 0.2719285784833474
$arr[23,1] = 0.2702162189359427
$arr[23,2] = 1.327933243004676
$ws.Range("L2:N25").Value = $arr

Write-Output "Applied case with 380 kV data"